# edit.ps1
# Applies the "Updated cryptos list" data refresh to Sheet1 of cryptos.xlsx.
# For each changed row, price (D) and/or 1h volume change (E) are updated to
# the latest scraped values; rows 13/14 (WrappedEther / Chainlink) also swap
# rank position (name, link, price, change all updated).
#
# D-column prices are forced to Text (NumberFormat "@") before assignment so
# plain-looking numeric strings (e.g. "0.285", "1.21") are NOT silently
# reinterpreted as numbers by Excel's smart input parsing -- the source data
# keeps these as text cells. The style is then reset back to "Normal" so the
# text-format override doesn't leave a lingering NumberFormat on the cell
# (matching the original/target workbook, where these cells carry no
# explicit style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.122.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.75%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.790.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.08%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "222.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.02%  "

# Row 6
$ws.Range("E6").Value = "  -0.45%  "

# Row 7
$ws.Range("E7").Value = "  +0.09%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.38"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.56%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.285"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.46%  "

# Row 10
$ws.Range("E10").Value = "  +0.41%  "

# Row 11
$ws.Range("E11").Value = "  +0.19%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.046.54"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.05%  "

# Row 13
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.98"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.68%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.788.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.13%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.626"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.17%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.085.73"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.95%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.90%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.09"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.52%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.11%  "

# Row 20
$ws.Range("E20").Value = "  -2.58%  "

# Row 22
$ws.Range("E22").Value = "  -1.28%  "

# Row 23
$ws.Range("E23").Value = "  -4.20%  "

# Row 24
$ws.Range("E24").Value = "  -1.21%  "

# Row 25
$ws.Range("E25").Value = "  -1.49%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.00%  "

# Row 27
$ws.Range("E27").Value = "  -1.84%  "

# Row 28
$ws.Range("E28").Value = "  -2.38%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0520"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.08%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.27%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.68"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.66%  "

# Row 33
$ws.Range("E33").Value = "  -3.60%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.80"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.10%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.393.30"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.79%  "

# Row 36
$ws.Range("E36").Value = "  +1.29%  "

# Row 37
$ws.Range("E37").Value = "  -1.35%  "

# Row 38
$ws.Range("E38").Value = "  -4.08%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "79.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.67%  "

# Row 40
$ws.Range("E40").Value = "  +0.82%  "

# Row 41
$ws.Range("E41").Value = "  -4.87%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.77%  "

# Row 43
$ws.Range("E43").Value = "  +1.43%  "

# Row 44
$ws.Range("E44").Value = "  -2.50%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0497"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.17%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "107.69"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.33%  "

# Row 47
$ws.Range("E47").Value = "  -0.99%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.947.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.63%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "12.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.08%  "

# Row 51
$ws.Range("E51").Value = "  +1.69%  "

